$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) and first column (A2:A6) labels: strip the " Diff-in-Diff" suffix
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "U"
$ws.Range("D1").Value = '$\pi$'
$ws.Range("E1").Value = "FFR"
$ws.Range("F1").Value = "A"

$ws.Range("A2").Value = "C"
$ws.Range("A3").Value = "U"
$ws.Range("A4").Value = '$\pi$'
$ws.Range("A5").Value = "FFR"
$ws.Range("A6").Value = "A"

# Helper: write a value that must remain TEXT even when it looks like a plain
# number (Excel would otherwise silently coerce "0.041" etc. into a numeric
# cell). Temporarily mark the cell as Text, assign, then restore the default
# "Normal" style so no stray number format lingers on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "C2" "0.839**"
Set-TextValue "D2" "-0.054**"
Set-TextValue "E2" "0.041"
Set-TextValue "F2" "0.079***"

# Row 3
Set-TextValue "B3" "0.084**"
Set-TextValue "D3" "-0.016**"
Set-TextValue "E3" "0.003"
Set-TextValue "F3" "-0.008*"

# Row 4
Set-TextValue "B4" "-1.158**"
Set-TextValue "C4" "-3.396**"
Set-TextValue "E4" "-0.11"
Set-TextValue "F4" "-0.231***"

# Row 5
Set-TextValue "B5" "0.633"
Set-TextValue "C5" "0.52"
Set-TextValue "D5" "-0.077"
Set-TextValue "F5" "-0.292***"

# Row 6
Set-TextValue "B6" "4.359***"
Set-TextValue "C6" "-4.271*"
Set-TextValue "D6" "-0.59***"
Set-TextValue "E6" "-1.054***"

# Row 7 (Constant)
Set-TextValue "B7" "-0.182"
Set-TextValue "C7" "-0.214"
Set-TextValue "D7" "0.006"
Set-TextValue "E7" "0.204**"
Set-TextValue "F7" "0.061"

# Row 8 (r2_adj) - genuine numeric values
$ws.Range("B8").Value = 0.71
$ws.Range("C8").Value = 0.19
$ws.Range("D8").Value = 0.6
$ws.Range("E8").Value = 0.43
$ws.Range("F8").Value = 0.79
